$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.013.02'
$ws.Range("E2").Value = '  -1.10%  '
$ws.Range("D3").Value = '2.476.25'
$ws.Range("E3").Value = '  -2.26%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.21'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.52%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  -2.22%  '
$ws.Range("D9").Value = '2.475.24'
$ws.Range("E9").Value = '  -2.22%  '
$ws.Range("E10").Value = '  -1.16%  '
$ws.Range("E12").Value = '  -2.23%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.329'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.19%  '
$ws.Range("D14").Value = '2.919.90'
$ws.Range("E14").Value = '  -2.31%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.25'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.67%  '
$ws.Range("D16").Value = '66.914.91'
$ws.Range("E16").Value = '  -1.01%  '
$ws.Range("E17").Value = '  -3.86%  '
$ws.Range("D18").Value = '2.459.05'
$ws.Range("E18").Value = '  -4.61%  '
$ws.Range("E19").Value = '  -7.18%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.35'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -9.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '349.77'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.24%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.01'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.69%  '
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.44'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.86%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.20'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -7.23%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.79'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.60%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.12'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.70%  '
$ws.Range("E28").Value = '  -26.13%  '
$ws.Range("D29").Value = '2.587.71'
$ws.Range("E29").Value = '  -2.80%  '
$ws.Range("E30").Value = '  -5.38%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '508.86'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.77%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.60'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.76%  '
$ws.Range("E33").Value = '  -5.00%  '
$ws.Range("E34").Value = '  -5.08%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '158.64'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.48%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.115'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -10.64%  '
$ws.Range("E39").Value = '  -6.37%  '
$ws.Range("E40").Value = '  -7.84%  '
$ws.Range("E41").Value = '  -0.32%  '
$ws.Range("E42").Value = '  -5.52%  '
$ws.Range("E43").Value = '  -5.33%  '
$ws.Range("E44").Value = '  -5.37%  '
$ws.Range("E45").Value = '  -4.16%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '38.77'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.57%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '141.01'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.03%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.43'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -7.26%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.511'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -7.02%  '
$ws.Range("D50").Value = '0.0₆0250'
$ws.Range("E50").Value = '  -8.34%  '
$ws.Range("E51").Value = '  -1.98%  '
